$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

# A21 / O21 / Q21 look like numbers or dates, but the source data stores
# them as plain text. Use the quote-prefix trick to force text entry, then
# clear the style back to Normal so no stray style index is left behind.
$ws.Cells.Item($row, 1).Value = "'2033555832"
$ws.Cells.Item($row, 1).Style = "Normal"

# B21 is present in the source but empty (an empty *text* cell). Writing ""
# directly gets pruned entirely by the engine, so use the quote-prefix
# trick (forces text type with empty content) then reset the style so it
# ends up a plain, un-styled empty text cell.
$ws.Cells.Item($row, 2).Value = "'"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "ddd"
$ws.Cells.Item($row, 4).Value = "Unité"
$ws.Cells.Item($row, 5).Value = "Unité"
$ws.Cells.Item($row, 6).Value = 10
$ws.Cells.Item($row, 7).Value = 100
$ws.Cells.Item($row, 8).Value = "Site principal"
$ws.Cells.Item($row, 9).Value = "E2"
$ws.Cells.Item($row, 10).Value = "E2"
$ws.Cells.Item($row, 11).Value = "FournX"
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = "Général"
$ws.Cells.Item($row, 14).Value = "Général"

$ws.Cells.Item($row, 15).Value = "'2033555832"
$ws.Cells.Item($row, 15).Style = "Normal"

$ws.Cells.Item($row, 16).Value = 0

$ws.Cells.Item($row, 17).Value = "'2025-06-04"
$ws.Cells.Item($row, 17).Style = "Normal"
